$d = $word.ActiveDocument

# --- Paragraph 1: fill in the existing empty trailing list paragraph ---
$n = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($n)
$r = $p1.Range
# The paragraph is empty, so setting LanguageID here stamps the
# paragraph-mark run properties (w:pPr/w:rPr/w:lang).
$r.LanguageID = "en-GB"

$r.InsertAfter("The buffer is empty at the end")
$r.LanguageID = "en-GB"

$r2 = $d.Range($r.End, $r.End)
$r2.InsertAfter(".")
$r2.LanguageID = "en-GB"
# Toggle a property to stop this run being silently merged with the
# previous, identically-formatted one.
$r2.Bold = 1
$r2.Bold = 0

# --- Paragraph 2: brand-new list paragraph after it ---
$r.InsertParagraphAfter()
$n2 = $d.Paragraphs.Count
$p2 = $d.Paragraphs.Item($n2)
$r3 = $p2.Range
# Again empty at this point, so this sets the new paragraph's mark rPr/lang.
$r3.LanguageID = "en-GB"

$r3.InsertAfter("When we leave the synchronisation, commands are only pushed and not pulled")
$r3.LanguageID = "en-GB"

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(" and the PushThread never end")
$r4.LanguageID = "en-GB"
$r4.Bold = 1
$r4.Bold = 0

$r5 = $d.Range($r4.End, $r4.End)
$r5.InsertAfter("s")
$r5.LanguageID = "en-GB"
$r5.Bold = 1
$r5.Bold = 0

$r6 = $d.Range($r5.End, $r5.End)
$r6.InsertAfter(".")
$r6.LanguageID = "en-GB"
$r6.Bold = 1
$r6.Bold = 0

Write-Host "Done"
